$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.573.48'
$ws.Range('E2').Value = '  +2.01%  '
$ws.Range('D3').Value = '2.944.90'
$ws.Range('E3').Value = '  +1.07%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''588.30'
$ws.Range('D6').Value = '''146.57'
$ws.Range('E6').Value = '  +3.09%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '2.945.73'
$ws.Range('E8').Value = '  +1.13%  '
$ws.Range('D9').Value = '''0.504'
$ws.Range('E9').Value = '  +1.33%  '
$ws.Range('D10').Value = '''6.91'
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('D11').Value = '''0.148'
$ws.Range('E11').Value = '  +5.70%  '
$ws.Range('D12').Value = '''0.434'
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('D13').Value = '''0.0000231'
$ws.Range('E13').Value = '  +4.62%  '
$ws.Range('D14').Value = '''32.04'
$ws.Range('E14').Value = '  -2.20%  '
$ws.Range('D15').Value = '''0.125'
$ws.Range('D16').Value = '3.431.36'
$ws.Range('E16').Value = '  +1.02%  '
$ws.Range('D17').Value = '62.560.26'
$ws.Range('E17').Value = '  +1.98%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.946.65'
$ws.Range('E18').Value = '  +1.19%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').Value = '''6.62'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('D20').Value = '''434.88'
$ws.Range('E20').Value = '  +1.25%  '
$ws.Range('D21').Value = '''13.35'
$ws.Range('E21').Value = '  -1.21%  '
$ws.Range('D22').Value = '''0.658'
$ws.Range('E22').Value = '  -0.85%  '
$ws.Range('D23').Value = '''6.93'
$ws.Range('E23').Value = '  -1.34%  '
$ws.Range('D24').Value = '''11.16'
$ws.Range('E24').Value = '  +3.20%  '
$ws.Range('D25').Value = '''80.06'
$ws.Range('E25').Value = '  -1.18%  '
$ws.Range('D26').Value = '''11.78'
$ws.Range('E26').Value = '  +1.81%  '
$ws.Range('D27').Value = '''2.09'
$ws.Range('E27').Value = '  -1.53%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').Value = '''7.17'
$ws.Range('E29').Value = '  +4.65%  '
$ws.Range('D30').Value = '''2.17'
$ws.Range('E30').Value = '  +1.35%  '
$ws.Range('D31').Value = '''2.58'
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('D32').Value = '0.0₃0991'
$ws.Range('E32').Value = '  +14.31%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '''26.11'
$ws.Range('E33').Value = '  -1.41%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '''0.107'
$ws.Range('E34').Value = '  +0.53%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').Value = '''0.990'
$ws.Range('E36').Value = '  -0.70%  '
$ws.Range('D37').Value = '''5.53'
$ws.Range('E37').Value = '  +0.25%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').Value = '''49.67'
$ws.Range('E38').Value = '  +0.35%  '
$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D39').Value = '''2.96'
$ws.Range('E39').Value = '  +1.67%  '
$ws.Range('D40').Value = '''1.98'
$ws.Range('E40').Value = '  +1.62%  '
$ws.Range('D41').Value = '''8.34'
$ws.Range('E41').Value = '  -0.53%  '
$ws.Range('D42').Value = '''0.114'
$ws.Range('E42').Value = '  -4.85%  '
$ws.Range('D43').Value = '''0.273'
$ws.Range('E43').Value = '  +0.99%  '
$ws.Range('D44').Value = '''39.45'
$ws.Range('E44').Value = '  -4.59%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.676.75'
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').Value = '''134.26'
$ws.Range('E46').Value = '  +0.64%  '
$ws.Range('D47').Value = '''0.0333'
$ws.Range('E47').Value = '  -1.47%  '
$ws.Range('D48').Value = '''351.23'
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('D50').Value = '''0.103'
$ws.Range('E50').Value = '  -0.46%  '
$ws.Range('D51').Value = '''22.36'
$ws.Range('E51').Value = '  -3.30%  '
